# Actualización automática 2025-09-29 14:30:10
$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H20").Value = 102.6
$ws1.Range("M20").Value = 3075.4
$ws1.Range("H35").Value = "2 de 33"

# --- Hoja "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F20").Value = 4007.38
$ws2.Range("F35").Value = 24347.64

# --- Hoja "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D6").Value = 908.1
$ws3.Range("E6").Value = -93.97656919112706
$ws3.Range("F6").Value = 1.115432827056404

$ws3.Range("D12").Value = 16743.84
$ws3.Range("E12").Value = 5689.915375176599
$ws3.Range("F12").Value = 0.7463681278493121

$ws3.Range("D15").Value = 24604.87
$ws3.Range("E15").Value = 14138.14881339592
$ws3.Range("F15").Value = 0.6350788026743163
